# Apply updated NATMI values per Dr Hou's advice (recomputed with 3 expressing cells
# instead of 1, for the Sema3a-Plxna2 ligand-receptor pair).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cols = @("E","G","H","I","J","K","M","N","O","P","Q","R","S","T")

# Row => ordered value list aligned with $cols
$rowData = @{
    2 = @(3, 0.8029023333333333, 2.408707, 0.06206726394886004, 0.06206726394886004, 3, 30.789758, 92.36927399999999, 0.5225738279799484, 0.5225738279799484, 24.72116854096866, 222.4905168687179, 0.03243472771399764, 0.03243472771399764)
    3 = @(3, 0.8029023333333333, 2.408707, 0.06206726394886004, 0.06206726394886004, 3, 2.247453666666667, 6.742361000000001, 0.03814451759567488, 0.03814451759567488, 1.804485793025222, 16.240372137227, 0.002367525841812689, 0.002367525841812689)
    4 = @(3, 0.8029023333333333, 2.408707, 0.06206726394886004, 0.06206726394886004, 3, 0.8899283333333333, 2.669785, 0.01510415430279821, 0.01510415430279821, 0.7145255353327777, 6.430729817994999, 0.0009374735318360868, 0.0009374735318360868)
    5 = @(3, 0.8029023333333333, 2.408707, 0.06206726394886004, 0.06206726394886004, 3, 1.615977666666667, 4.847932999999999, 0.02742690069860586, 0.02742690069860585, 1.297472239181222, 11.677250152631, 0.001702312684959544, 0.001702312684959543)
    6 = @(3, 0.8029023333333333, 2.408707, 0.06206726394886004, 0.06206726394886004, 3, 0.2964566666666666, 0.88937, 0.005031559362375489, 0.005031559362375489, 0.2380257493988888, 2.14223174459, 0.0003122951230189174, 0.0003122951230189174)
    7 = @(3, 0.8029023333333333, 2.408707, 0.06206726394886004, 0.06206726394886004, 3, 23.079867, 69.23960099999999, 0.3917190400605971, 0.3917190400605971, 18.530879067323, 166.777911605907, 0.02431292905323516, 0.02431292905323516)
    8 = @(3, 2.159929333333333, 6.479788, 0.1669703754456877, 0.1669703754456877, 3, 30.789758, 92.36927399999999, 0.5225738279799484, 0.5225738279799484, 66.50370147043466, 598.5333132339119, 0.08725434825590218, 0.08725434825590218)
    9 = @(3, 2.159929333333333, 6.479788, 0.1669703754456877, 0.1669703754456877, 3, 2.247453666666667, 6.742361000000001, 0.03814451759567488, 0.03814451759567488, 4.854341099940889, 43.689069899468, 0.006369004424144474, 0.006369004424144474)
    10 = @(3, 2.159929333333333, 6.479788, 0.1669703754456877, 0.1669703754456877, 3, 0.8899283333333333, 2.669785, 0.01510415430279821, 0.01510415430279821, 1.922182311731111, 17.29964080558, 0.002521946314727816, 0.002521946314727816)
    11 = @(3, 2.159929333333333, 6.479788, 0.1669703754456877, 0.1669703754456877, 3, 1.615977666666667, 4.847932999999999, 0.02742690069860586, 0.02742690069860585, 3.490397564244889, 31.413578078204, 0.004579479906957814, 0.004579479906957813)
    12 = @(3, 2.159929333333333, 6.479788, 0.1669703754456877, 0.1669703754456877, 3, 0.2964566666666666, 0.88937, 0.005031559362375489, 0.005031559362375489, 0.6403254503955556, 5.76292905356, 0.0008401213558131001, 0.0008401213558131001)
    13 = @(3, 2.159929333333333, 6.479788, 0.1669703754456877, 0.1669703754456877, 3, 23.079867, 69.23960099999999, 0.3917190400605971, 0.3917190400605971, 49.85088174273199, 448.657935684588, 0.06540547518814226, 0.06540547518814226)
    14 = @(3, 9.973171666666666, 29.919515, 0.7709623606054523, 0.7709623606054523, 3, 30.789758, 92.36927399999999, 0.5225738279799484, 0.5225738279799484, 307.0715421091232, 2763.64387898211, 0.4028847520100486, 0.4028847520100486)
    15 = @(3, 9.973171666666666, 29.919515, 0.7709623606054523, 0.7709623606054523, 3, 2.247453666666667, 6.742361000000001, 0.03814451759567488, 0.03814451759567488, 22.41424123054611, 201.728171074915, 0.02940798732971772, 0.02940798732971772)
    16 = @(3, 9.973171666666666, 29.919515, 0.7709623606054523, 0.7709623606054523, 3, 0.8899283333333333, 2.669785, 0.01510415430279821, 0.01510415430279821, 8.875408039363888, 79.878672354275, 0.01164473445623431, 0.01164473445623431)
    17 = @(3, 9.973171666666666, 29.919515, 0.7709623606054523, 0.7709623606054523, 3, 1.615977666666667, 4.847932999999999, 0.02742690069860586, 0.02742690069860585, 16.11642267916611, 145.047804112495, 0.0211451081066885, 0.0211451081066885)
    18 = @(3, 9.973171666666666, 29.919515, 0.7709623606054523, 0.7709623606054523, 3, 0.2964566666666666, 0.88937, 0.005031559362375489, 0.005031559362375489, 2.956613228394444, 26.60951905555, 0.003879142883543471, 0.003879142883543471)
    19 = @(3, 9.973171666666666, 29.919515, 0.7709623606054523, 0.7709623606054523, 3, 23.079867, 69.23960099999999, 0.3917190400605971, 0.3917190400605971, 230.1794756348349, 2071.615280713514, 0.3020006358192197, 0.3020006358192197)
}

foreach ($row in $rowData.Keys) {
    $values = $rowData[$row]
    for ($i = 0; $i -lt $cols.Length; $i++) {
        $ws.Range($cols[$i] + $row).Value = $values[$i]
    }
}

Write-Output "Updated $($rowData.Keys.Count) rows"